# Update countries & provincias Spain
# - Swap the displayed country names for rows 110/111 (Montenegro <-> Vietnam)
# - Update "Datos actualizados" timestamp in A1
# - Update several countries' statistics (row 24 India, row 33 Rumania,
#   row 45 Finlandia, row 68 Lituania, row 110/111 Montenegro/Vietnam,
#   row 129 El Salvador, row 152 Zambia)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: updated timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 13:52"

# --- Row 24: India ---
$ws.Range("B24").Value = 6725
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 641
$ws.Range("E24").Value = 5855
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 229

# --- Row 33: Rumania ---
$ws.Range("E33").Value = 4473
$ws.Range("G33").Value = 17
$ws.Range("H33").Value = 265

# --- Row 45: Finlandia ---
$ws.Range("E45").Value = 2421
$ws.Range("G45").Value = 6
$ws.Range("H45").Value = 48

# --- Row 68: Lituania ---
$ws.Range("D68").Value = 54
$ws.Range("E68").Value = 923
$ws.Range("F68").Value = 14
$ws.Range("G68").Value = 6
$ws.Range("H68").Value = 22

# --- Rows 110/111: Montenegro and Vietnam swap places (Vietnam now listed
#     first, with refreshed stats; Montenegro drops to the row below,
#     keeping its previous stats) ---
$ws.Range("A110").Value = "Vietnam"
$ws.Range("B110").Value = 257
$ws.Range("C110").Value = 2
$ws.Range("D110").Value = 144
$ws.Range("E110").Value = 113
$ws.Range("F110").Value = 8
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 0

$ws.Range("A111").Value = "Montenegro"
$ws.Range("B111").Value = 255
$ws.Range("C111").Value = 3
$ws.Range("D111").Value = 4
$ws.Range("E111").Value = 249
$ws.Range("F111").Value = 7
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 2

# --- Row 129: El Salvador ---
$ws.Range("F129").Value = 4

# --- Row 152: Zambia ---
$ws.Range("E152").Value = 13
$ws.Range("G152").Value = 1
$ws.Range("H152").Value = 2
